# "Generate Report for Handback" — mark the zh-cn and de-de localization
# rows as handed back: update status, populate the target/handback file
# columns (with a hyperlink on the target file, matching the source-file
# hyperlink), stamp the handback datetime, and widen the columns that now
# hold longer filenames.

$wb = $excel.ActiveWorkbook

$mdFile     = "d0ca1dd6-2ade-4e81-a460-1e79887cd4d5.md"
$mdUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/d0ca1dd6-2ade-4e81-a460-1e79887cd4d5.md"
$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# zh-cn: status, latest target file (hyperlinked), latest handback file,
# latest handback datetime
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFile)

$wsZhCn.Range("J2").Value = "d0ca1dd6-2ade-4e81-a460-1e79887cd4d5.4f97ae72bdc331ac2cf7ad237d892bf152f05406.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-22 17:00:57"

# ---------------------------------------------------------------------
# de-de: status, latest target file (hyperlinked), latest handback file,
# latest handback datetime
# ---------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFile)

$wsDeDe.Range("J2").Value = "d0ca1dd6-2ade-4e81-a460-1e79887cd4d5.4f97ae72bdc331ac2cf7ad237d892bf152f05406.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-22 17:01:26"

# ---------------------------------------------------------------------
# Widen columns to fit the newly-populated long file names
# ---------------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I1").ColumnWidth = 40
$wsZhCn.Range("J1").ColumnWidth = 40

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I1").ColumnWidth = 40
$wsDeDe.Range("J1").ColumnWidth = 40

Write-Host "Handback report generated for zh-cn and de-de."
